$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" values under duplicate_image_filename (column E) for rows 2 through 21
$ws.Range("E2:E21").Value = "NA"

# F1 is an empty-string placeholder cell in the original file; re-clearing it
# keeps it an empty cell (avoids an unrelated save-time artifact turning it
# into a non-blank value, which is not part of this change).
$ws.Range("F1").ClearContents()
